# Adds three new market sheets (Netherlands, Austria, Denmark) to the workbook,
# copied from the existing "Italy" sheet template (same layout/column widths as
# the other per-country sheets), positioned after "Greece".
#
# Austria's sheet additionally drops the "PR1D2-Unmonitored" row (row 10),
# matching the reduced 11-row layout used for that market's test data.

$wb = $excel.ActiveWorkbook
$template = $wb.Worksheets.Item("Italy")
$greece = $wb.Worksheets.Item("Greece")

# --- Netherlands -----------------------------------------------------------
$template.Copy($null, $greece)
$netherlands = $wb.Worksheets.Item("Italy (2)")
$netherlands.Name = "Netherlands"
$netherlands.Range("B4").Value = "NGC-3144/T2176"
$netherlands.Range("B2").Value = "Netherlands Market"
$netherlands.Activate()
$netherlands.Range("B4").Select()

# --- Austria -----------------------------------------------------------
$template.Copy($null, $netherlands)
$austria = $wb.Worksheets.Item("Italy (2)")
$austria.Name = "Austria"
$austria.Range("B4").Value = "NGC-3817/"
$austria.Range("B2").Value = "Austria Market"
$austria.Rows.Item(10).Delete()
$austria.Activate()
$austria.Range("C14").Select()

# --- Denmark -----------------------------------------------------------
$template.Copy($null, $austria)
$denmark = $wb.Worksheets.Item("Italy (2)")
$denmark.Name = "Denmark"
$denmark.Range("B4").Value = "NGC-2913/"
$denmark.Range("B2").Value = "Denmark Market"
$denmark.Activate()
$denmark.Range("B4").Select()

# Austria ends up the active/selected tab.
$austria.Activate()
